# Update "想去人数" (want-to-go count) figures for three events that
# appear on both the "展览" sheet and the "全部类型" sheet.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions): rows 2, 4, 5 hold the affected events.
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 325
$wsExhibition.Range("F4").Value = 1319
$wsExhibition.Range("F5").Value = 643

# Sheet "全部类型" (All types): same events, but row 5 holds an extra
# concert entry, so the ANE·DACG event sits in row 6 instead of row 5.
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 325
$wsAll.Range("F4").Value = 1319
$wsAll.Range("F6").Value = 643
